# Applies the "Tests for _collects_nutrient_ratio_conflicts() all passing!" edit
# to the flag conflict logic table (docs/flag_conflict_logic_table~v2.xlsx).
#
# Summary of content changes on Sheet1:
#   - New "Direct Alias" Flag Type category (rows 9, 11, 13, 14 in col A).
#   - H9 reworded ("," -> " or ").
#   - G15 changed from "One or More" to "Two or More".
#   - Test names filled in / renamed in column I for rows 9-15.
#   - Removed the now-unused conditional formatting rules on column A that
#     highlighted "Has DOF" / "Direct Alias".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = "Direct Alias"
$ws.Range("H9").Value = "Related nutrient in 'need_zero' or 'need_non_zero'  to oppose implication."

# Row 10
$ws.Range("I10").Value = "test_none_to_false_with_single_nutrient_correctly_categorises_opposing_implication()"

# Row 11
$ws.Range("A11").Value = "Direct Alias"
$ws.Range("I11").Value = "test_true_to_none_with_direct_alias_single_nutrient_correctly_categorises_need_undefining()"

# Row 9 (test name, renamed after rows 10-11 were filled in)
$ws.Range("I9").Value = "test_true_to_false_with_direct_alias_single_nutrient_correctly_categorises_opposing_implication()"

# Row 12
$ws.Range("I12").Value = "test_false_to_none_collects_all_defined_opposing_nutrients_in_need_undefining()"

# Row 13
$ws.Range("A13").Value = "Direct Alias"
$ws.Range("I13").Value = "test_true_to_false_with_direct_alias_multiple_related_nutrients_collects_all_in_preventing_flag_false()"

# Row 14
$ws.Range("A14").Value = "Direct Alias"
$ws.Range("I14").Value = "test_true_to_none_with_direct_alias_multiple_related_nutrients_collects_all_in_preventing_flag_undefine()"

# Row 15
$ws.Range("G15").Value = "Two or More"
$ws.Range("I15").Value = "test_none_to_false_with_multiple_undefined_nutrients_collects_all_in_preventing_flag_false()"

# Remove the "Has DOF " / "Direct Alias" conditional formatting on column A
# (its backing dxfs become unused, as in the target workbook).
$ws.Range("A1:A1048576").FormatConditions.Delete()

# Leave the active selection on I15, matching where the user last worked.
$ws.Range("I15").Select()
